# Refresh the "cryptos" price/volume snapshot (GitHub Actions style update).
# For every "Price" (column D) cell whose new value would otherwise be
# auto-coerced into a floating point Number by Excel (losing exact text
# such as trailing zeros, e.g. "1.00" -> 1), the NumberFormat is forced to
# "@" (Text) immediately before the value is written so the stored cell
# keeps the same text type/content as the original inline string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.341.37'
$ws.Range('E2').Value = '  -1.42%  '
$ws.Range('D3').Value = '1.591.76'
$ws.Range('E3').Value = '  -0.45%  '
$ws.Range('E4').Value = '  -0.48%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '210.44'
$ws.Range('E5').Value = '  -0.58%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.504'
$ws.Range('E6').Value = '  -1.69%  '
$ws.Range('E7').Value = '  -0.46%  '
$ws.Range('E8').Value = '  -1.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.246'
$ws.Range('E9').Value = '  -0.40%  '
$ws.Range('E10').Value = '  +0.20%  '
$ws.Range('E11').Value = '  -0.61%  '
$ws.Range('E13').Value = '  +0.45%  '
$ws.Range('D14').Value = '1.574.13'
$ws.Range('E14').Value = '  -1.52%  '
$ws.Range('E15').Value = '  -1.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.68'
$ws.Range('E16').Value = '  -0.55%  '
$ws.Range('D17').Value = '26.342.10'
$ws.Range('E17').Value = '  -1.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.47'
$ws.Range('E19').Value = '  +4.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '211.86'
$ws.Range('E20').Value = '  +1.70%  '
$ws.Range('E21').Value = '  -0.50%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.28'
$ws.Range('E22').Value = '  -0.25%  '
$ws.Range('E23').Value = '  -1.67%  '
$ws.Range('E24').Value = '  -1.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.94'
$ws.Range('E25').Value = '  +0.86%  '
$ws.Range('E26').Value = '  -0.43%  '
$ws.Range('E27').Value = '  -0.78%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.113'
$ws.Range('E28').Value = '  -1.18%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.29'
$ws.Range('E29').Value = '  -0.19%  '
$ws.Range('E30').Value = '  -0.30%  '
$ws.Range('E31').Value = '  -1.03%  '
$ws.Range('E32').Value = '  -0.78%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.99'
$ws.Range('E33').Value = '  +0.44%  '
$ws.Range('D34').Value = '1.297.03'
$ws.Range('E34').Value = '  +1.72%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.614'
$ws.Range('E35').Value = '  +3.99%  '
$ws.Range('E36').Value = '  -1.69%  '
$ws.Range('E37').Value = '  -0.94%  '
$ws.Range('E38').Value = '  -0.75%  '
$ws.Range('E39').Value = '  -12.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.816'
$ws.Range('E40').Value = '  -0.92%  '
$ws.Range('E41').Value = '  -0.51%  '
$ws.Range('E42').Value = '  +2.89%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.14'
$ws.Range('E43').Value = '  -2.34%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '62.68'
$ws.Range('E44').Value = '  +0.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.761'
$ws.Range('E45').Value = '  -1.82%  '
$ws.Range('D46').Value = '1.726.90'
$ws.Range('E46').Value = '  -0.46%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '88.43'
$ws.Range('E47').Value = '  -2.17%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.51'
$ws.Range('E48').Value = '  -3.64%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0998'
$ws.Range('E49').Value = '  -2.33%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0505'
$ws.Range('E50').Value = '  -1.47%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').Value = '  -0.42%  '

Write-Host "Applied cell changes"
